$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.794.31'
$ws.Range("E2").Value = '  -2.65%  '

# Row 3
$ws.Range("D3").Value = '1.777.72'
$ws.Range("E3").Value = '  -2.98%  '

# Row 4
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.39%  '

# Row 5
$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").Value = '1.005'
$ws.Range("E5").Value = '  +0.39%  '

# Row 6
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '307.49'
$ws.Range("E6").Value = '  -1.77%  '

# Row 7
$ws.Range("D7").Value = '0.4430'
$ws.Range("E7").Value = '  +3.18%  '

# Row 8
$ws.Range("D8").Value = '0.3637'
$ws.Range("E8").Value = '  -0.86%  '

# Row 9
$ws.Range("D9").Value = '0.07205'
$ws.Range("E9").Value = '  -1.07%  '

# Row 10
$ws.Range("D10").Value = '0.8378'
$ws.Range("E10").Value = '  -3.56%  '

# Row 11
$ws.Range("D11").Value = '20.31'
$ws.Range("E11").Value = '  -1.88%  '

# Row 12
$ws.Range("D12").Value = '1.843.46'
$ws.Range("E12").Value = '  -1.77%  '

# Row 13
$ws.Range("D13").Value = '5.271'
$ws.Range("E13").Value = '  -2.59%  '

# Row 14
$ws.Range("D14").Value = '6.360'
$ws.Range("E14").Value = '  -2.84%  '

# Row 15
$ws.Range("D15").Value = '0.06802'
$ws.Range("E15").Value = '  -1.99%  '

# Row 16
$ws.Range("E16").Value = '  +0.68%  '

# Row 17
$ws.Range("D17").Value = '79.67'
$ws.Range("E17").Value = '  -1.35%  '

# Row 18
$ws.Range("D18").Value = '0.000008704'
$ws.Range("E18").Value = '  -2.52%  '

# Row 19
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.48%  '

# Row 20
$ws.Range("D20").Value = '14.99'
$ws.Range("E20").Value = '  -2.90%  '

# Row 21
$ws.Range("D21").Value = '26.892.78'
$ws.Range("E21").Value = '  -2.37%  '

# Row 22
$ws.Range("D22").Value = '5.040'
$ws.Range("E22").Value = '  -2.15%  '

# Row 23
$ws.Range("D23").Value = '11.09'
$ws.Range("E23").Value = '  +1.69%  '

# Row 24
$ws.Range("D24").Value = '2.033.20'
$ws.Range("E24").Value = '  -3.10%  '

# Row 25
$ws.Range("D25").Value = '1.918'
$ws.Range("E25").Value = '  -3.15%  '

# Row 26
$ws.Range("D26").Value = '153.61'
$ws.Range("E26").Value = '  -0.68%  '

# Row 27
$ws.Range("D27").Value = '18.21'
$ws.Range("E27").Value = '  -3.64%  '

# Row 28
$ws.Range("D28").Value = '115.31'
$ws.Range("E28").Value = '  +0.96%  '

# Row 29
$ws.Range("D29").Value = '5.046'
$ws.Range("E29").Value = '  -1.98%  '

# Row 30
$ws.Range("D30").Value = '1.635'
$ws.Range("E30").Value = '  -11.27%  '

# Row 31
$ws.Range("D31").Value = '0.09051'
$ws.Range("E31").Value = '  +2.23%  '

# Row 32
$ws.Range("D32").Value = '0.7247'
$ws.Range("E32").Value = '  -4.17%  '

# Row 33
$ws.Range("D33").Value = '2.842'
$ws.Range("E33").Value = '  -4.54%  '

# Row 34
$ws.Range("D34").Value = '4.336'
$ws.Range("E34").Value = '  -4.78%  '

# Row 35
$ws.Range("D35").Value = '1.095'
$ws.Range("E35").Value = '  -3.76%  '

# Row 36
$ws.Range("D36").Value = '1.005'
$ws.Range("E36").Value = '  +0.45%  '

# Row 37
$ws.Range("D37").Value = '1.076'
$ws.Range("E37").Value = '  -1.19%  '

# Row 38
$ws.Range("D38").Value = '0.01891'
$ws.Range("E38").Value = '  -2.61%  '

# Row 39
$ws.Range("D39").Value = '0.05096'
$ws.Range("E39").Value = '  -4.56%  '

# Row 40
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.1611'
$ws.Range("E40").Value = '  -3.46%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.4919'
$ws.Range("E41").Value = '  -3.54%  '

# Row 42
$ws.Range("D42").Value = '2.569'
$ws.Range("E42").Value = '  -8.23%  '

# Row 43
$ws.Range("D43").Value = '6.117'
$ws.Range("E43").Value = '  -7.20%  '

# Row 44
$ws.Range("D44").Value = '7.929'
$ws.Range("E44").Value = '  -5.48%  '

# Row 45
$ws.Range("D45").Value = '104.83'
$ws.Range("E45").Value = '  -1.44%  '

# Row 46
$ws.Range("D46").Value = '1.005'
$ws.Range("E46").Value = '  +0.47%  '

# Row 47
$ws.Range("D47").Value = '10.04'
$ws.Range("E47").Value = '  -3.96%  '

# Row 48
$ws.Range("D48").Value = '0.06230'
$ws.Range("E48").Value = '  -4.23%  '

# Row 49
$ws.Range("D49").Value = '0.4493'
$ws.Range("E49").Value = '  -4.49%  '

# Row 50
$ws.Range("D50").Value = '1.576'
$ws.Range("E50").Value = '  -2.91%  '

# Row 51
$ws.Range("D51").Value = '1.732'
$ws.Range("E51").Value = '  -0.53%  '
